$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 6.654043666666666
$ws.Range("H2").Value = 19.962131
$ws.Range("I2").Value = 0.3091924566209486
$ws.Range("J2").Value = 0.3091924566209486
$ws.Range("M2").Value = 43.68636333333333
$ws.Range("N2").Value = 131.05909
$ws.Range("O2").Value = 0.4655970190670373
$ws.Range("P2").Value = 0.4655970190670373
$ws.Range("Q2").Value = 290.6909692578656
$ws.Range("R2").Value = 2616.21872332079
$ws.Range("S2").Value = 0.1439590861207279
$ws.Range("T2").Value = 0.1439590861207279
$ws.Range("G3").Value = 6.654043666666666
$ws.Range("H3").Value = 19.962131
$ws.Range("I3").Value = 0.3091924566209486
$ws.Range("J3").Value = 0.3091924566209486
$ws.Range("O3").Value = 0.01284473362738172
$ws.Range("P3").Value = 0.01284473362738172
$ws.Range("Q3").Value = 8.019484479270444
$ws.Range("R3").Value = 72.175360313434
$ws.Range("S3").Value = 0.003971494744891862
$ws.Range("T3").Value = 0.003971494744891862
$ws.Range("G4").Value = 6.654043666666666
$ws.Range("H4").Value = 19.962131
$ws.Range("I4").Value = 0.3091924566209486
$ws.Range("J4").Value = 0.3091924566209486
$ws.Range("M4").Value = 4.981224333333333
$ws.Range("N4").Value = 14.943673
$ws.Range("O4").Value = 0.05308849315764798
$ws.Range("P4").Value = 0.05308849315764799
$ws.Range("Q4").Value = 33.14528422746255
$ws.Range("R4").Value = 298.307558047163
$ws.Range("S4").Value = 0.0164145616177176
$ws.Range("T4").Value = 0.0164145616177176
$ws.Range("G5").Value = 6.654043666666666
$ws.Range("H5").Value = 19.962131
$ws.Range("I5").Value = 0.3091924566209486
$ws.Range("J5").Value = 0.3091924566209486
$ws.Range("M5").Value = 43.95590833333333
$ws.Range("N5").Value = 131.867725
$ws.Range("O5").Value = 0.4684697541479331
$ws.Range("P5").Value = 0.4684697541479331
$ws.Range("Q5").Value = 292.4845334579972
$ws.Range("R5").Value = 2632.360801121975
$ws.Range("S5").Value = 0.1448473141376112
$ws.Range("T5").Value = 0.1448473141376113
$ws.Range("I6").Value = 0.09233579784218476
$ws.Range("J6").Value = 0.09233579784218476
$ws.Range("M6").Value = 43.68636333333333
$ws.Range("N6").Value = 131.05909
$ws.Range("O6").Value = 0.4655970190670373
$ws.Range("P6").Value = 0.4655970190670373
$ws.Range("Q6").Value = 86.81059966753556
$ws.Range("R6").Value = 781.2953970078199
$ws.Range("S6").Value = 0.0429912722284978
$ws.Range("T6").Value = 0.0429912722284978
$ws.Range("I7").Value = 0.09233579784218476
$ws.Range("J7").Value = 0.09233579784218476
$ws.Range("O7").Value = 0.01284473362738172
$ws.Range("P7").Value = 0.01284473362738172
$ws.Range("S7").Value = 0.001186028727554631
$ws.Range("T7").Value = 0.001186028727554631
$ws.Range("I8").Value = 0.09233579784218476
$ws.Range("J8").Value = 0.09233579784218476
$ws.Range("M8").Value = 4.981224333333333
$ws.Range("N8").Value = 14.943673
$ws.Range("O8").Value = 0.05308849315764798
$ws.Range("P8").Value = 0.05308849315764799
$ws.Range("Q8").Value = 9.898353592761556
$ws.Range("R8").Value = 89.085182334854
$ws.Range("S8").Value = 0.004901968371950793
$ws.Range("T8").Value = 0.004901968371950794
$ws.Range("I9").Value = 0.09233579784218476
$ws.Range("J9").Value = 0.09233579784218476
$ws.Range("M9").Value = 43.95590833333333
$ws.Range("N9").Value = 131.867725
$ws.Range("O9").Value = 0.4684697541479331
$ws.Range("P9").Value = 0.4684697541479331
$ws.Range("Q9").Value = 87.34622134217223
$ws.Range("R9").Value = 786.11599207955
$ws.Range("S9").Value = 0.04325652851418155
$ws.Range("T9").Value = 0.04325652851418155
$ws.Range("G10").Value = 0.9593116666666667
$ws.Range("H10").Value = 2.877935
$ws.Range("I10").Value = 0.04457619242381535
$ws.Range("J10").Value = 0.04457619242381536
$ws.Range("M10").Value = 43.68636333333333
$ws.Range("N10").Value = 131.05909
$ws.Range("O10").Value = 0.4655970190670373
$ws.Range("P10").Value = 0.4655970190670373
$ws.Range("Q10").Value = 41.90883801990555
$ws.Range("R10").Value = 377.17954217915
$ws.Range("S10").Value = 0.02075454231388708
$ws.Range("T10").Value = 0.02075454231388708
$ws.Range("G11").Value = 0.9593116666666667
$ws.Range("H11").Value = 2.877935
$ws.Range("I11").Value = 0.04457619242381535
$ws.Range("J11").Value = 0.04457619242381536
$ws.Range("O11").Value = 0.01284473362738172
$ws.Range("P11").Value = 0.01284473362738172
$ws.Range("Q11").Value = 1.156166897454445
$ws.Range("R11").Value = 10.40550207709
$ws.Range("S11").Value = 0.0005725693178068193
$ws.Range("T11").Value = 0.0005725693178068194
$ws.Range("G12").Value = 0.9593116666666667
$ws.Range("H12").Value = 2.877935
$ws.Range("I12").Value = 0.04457619242381535
$ws.Range("J12").Value = 0.04457619242381536
$ws.Range("M12").Value = 4.981224333333333
$ws.Range("N12").Value = 14.943673
$ws.Range("O12").Value = 0.05308849315764798
$ws.Range("P12").Value = 0.05308849315764799
$ws.Range("Q12").Value = 4.778546617250556
$ws.Range("R12").Value = 43.006919555255
$ws.Range("S12").Value = 0.002366482886485721
$ws.Range("T12").Value = 0.002366482886485722
$ws.Range("G13").Value = 0.9593116666666667
$ws.Range("H13").Value = 2.877935
$ws.Range("I13").Value = 0.04457619242381535
$ws.Range("J13").Value = 0.04457619242381536
$ws.Range("M13").Value = 43.95590833333333
$ws.Range("N13").Value = 131.867725
$ws.Range("O13").Value = 0.4684697541479331
$ws.Range("P13").Value = 0.4684697541479331
$ws.Range("Q13").Value = 42.16741568309723
$ws.Range("R13").Value = 379.506741147875
$ws.Range("S13").Value = 0.02088259790563574
$ws.Range("T13").Value = 0.02088259790563574
$ws.Range("G14").Value = 11.92023
$ws.Range("H14").Value = 35.76069
$ws.Range("I14").Value = 0.5538955531130513
$ws.Range("J14").Value = 0.5538955531130514
$ws.Range("M14").Value = 43.68636333333333
$ws.Range("N14").Value = 131.05909
$ws.Range("O14").Value = 0.4655970190670373
$ws.Range("P14").Value = 0.4655970190670373
$ws.Range("Q14").Value = 520.7514987968999
$ws.Range("R14").Value = 4686.763489172099
$ws.Range("S14").Value = 0.2578921184039245
$ws.Range("T14").Value = 0.2578921184039246
$ws.Range("G15").Value = 11.92023
$ws.Range("H15").Value = 35.76069
$ws.Range("I15").Value = 0.5538955531130513
$ws.Range("J15").Value = 0.5538955531130514
$ws.Range("O15").Value = 0.01284473362738172
$ws.Range("P15").Value = 0.01284473362738172
$ws.Range("Q15").Value = 14.36631682374
$ws.Range("R15").Value = 129.29685141366
$ws.Range("S15").Value = 0.007114640837128407
$ws.Range("T15").Value = 0.007114640837128408
$ws.Range("G16").Value = 11.92023
$ws.Range("H16").Value = 35.76069
$ws.Range("I16").Value = 0.5538955531130513
$ws.Range("J16").Value = 0.5538955531130514
$ws.Range("M16").Value = 4.981224333333333
$ws.Range("N16").Value = 14.943673
$ws.Range("O16").Value = 0.05308849315764798
$ws.Range("P16").Value = 0.05308849315764799
$ws.Range("Q16").Value = 59.37733973492999
$ws.Range("R16").Value = 534.39605761437
$ws.Range("S16").Value = 0.02940548028149386
$ws.Range("T16").Value = 0.02940548028149387
$ws.Range("G17").Value = 11.92023
$ws.Range("H17").Value = 35.76069
$ws.Range("I17").Value = 0.5538955531130513
$ws.Range("J17").Value = 0.5538955531130514
$ws.Range("M17").Value = 43.95590833333333
$ws.Range("N17").Value = 131.867725
$ws.Range("O17").Value = 0.4684697541479331
$ws.Range("P17").Value = 0.4684697541479331
$ws.Range("Q17").Value = 523.96453719225
$ws.Range("R17").Value = 4715.680834730249
$ws.Range("S17").Value = 0.2594833135905045
$ws.Range("T17").Value = 0.2594833135905046
